$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantities / values per the diff
$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 2

$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 1

$ws.Range("B5").Value = 7
$ws.Range("C5").Value = 8

# Update the active selection to D3
$ws.Range("D3").Select() | Out-Null
